{"js": "// The resume's \"ABOUT\" paragraph reads:\n//   \"... Developer for over 2 years. As a mobile ...\"\n// The commit bumps the number of years of experience from 2 to 3, i.e.\n//   \" for over 2 \" -> \" for over 3 \"\n// (the author's tool happened to split that run into three runs with\n// identical formatting, but the visible/semantic text change is simply\n// the digit \"2\" -> \"3\" inside the existing \" for over 2 \" run).\n\nconst body = context.document.body;\n\n// Narrow the search to the exact phrase so we don't touch any other\n// standalone \"2\" elsewhere in the document.\nconst phraseResults = body.search(\" for over 2 \", { matchCase: true, matchWholeWord: false });\nphraseResults.load(\"text\");\nawait context.sync();\n\nif (phraseResults.items.length === 0) {\n  throw new Error('Could not find \" for over 2 \" in the document body.');\n}\n\nconst phraseRange = phraseResults.items[0];\n\n// Within that exact phrase, find the \"2\" itself and replace only that\n// character, preserving the surrounding \" for over \" / \" \" text and the\n// run's original formatting (bold Arial 12pt black, per the rPr in the\n// source XML).\nconst digitResults = phraseRange.search(\"2\", { matchCase: true });\nawait context.sync();\n\nif (digitResults.items.length === 0) {\n  throw new Error('Could not find \"2\" inside the \" for over 2 \" phrase.');\n}\n\ndigitResults.items[0].insertText(\"3\", \"Replace\");\nawait context.sync();\n", "ps1": "# The resume's \"ABOUT\" paragraph reads:\n#   \"... Developer for over 2 years. As a mobile ...\"\n# The commit bumps the number of years of experience from 2 to 3, i.e.\n#   \" for over 2 \" -> \" for over 3 \"\n# (the author's tool happened to split that run into three runs with\n# identical formatting, but the visible/semantic text change is simply\n# the digit \"2\" -> \"3\" inside the existing \" for over 2 \" run).\n\n$d = $word.ActiveDocument\n\n# Locate the exact phrase first so we never touch an unrelated standalone\n# \"2\" elsewhere in the document (e.g. dates like \"july/2023\").\n$phrase = $d.Content\n$phraseFound = $phrase.Find.Execute(\" for over 2 \", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n\nif (-not $phraseFound) {\n    throw 'Could not find \" for over 2 \" in the document.'\n}\n\n# Re-anchor a fresh Range over exactly that hit, then search for \"2\" within\n# it only (a Range.Find.Execute re-run on the SAME range would keep scanning\n# forward through the rest of the story, so scope a new Range first).\n$phraseStart = $phrase.Start\n$phraseEnd = $phrase.End\n$scoped = $d.Range($phraseStart, $phraseEnd)\n\n$digitFound = $scoped.Find.Execute(\"2\", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n\nif (-not $digitFound) {\n    throw 'Could not find \"2\" inside the \" for over 2 \" phrase.'\n}\n\n# Replace only the digit, preserving the run's original formatting\n# (bold Arial 12pt black) and the surrounding \" for over \" / \" \" text.\n$scoped.Text = \"3\"\n"}
